$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new value would otherwise be
# auto-detected as a plain number by Excel (so they stay text, as in the source).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"

# Apply the updated cell values from the crypto-price refresh.
$ws.Range("D2").Value = "64.396.22"
$ws.Range("E2").Value = "  +0.17%  "
$ws.Range("D3").Value = "3.500.09"
$ws.Range("E3").Value = "  +0.30%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "590.30"
$ws.Range("E5").Value = "  +0.56%  "
$ws.Range("D6").Value = "134.25"
$ws.Range("E6").Value = "  -0.10%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  +0.66%  "
$ws.Range("D9").Value = "7.64"
$ws.Range("E9").Value = "  +6.19%  "
$ws.Range("E10").Value = "  +1.01%  "
$ws.Range("D11").Value = "0.392"
$ws.Range("E11").Value = "  +4.28%  "
$ws.Range("D12").Value = "4.097.54"
$ws.Range("E12").Value = "  +0.35%  "
$ws.Range("E13").Value = "  +0.71%  "
$ws.Range("E14").Value = "  +0.58%  "
$ws.Range("D15").Value = "3.502.43"
$ws.Range("E15").Value = "  +0.30%  "
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "64.385.29"
$ws.Range("E16").Value = "  +0.10%  "
$ws.Range("B17").Value = "Avalanche"
$ws.Range("C17").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D17").Value = "25.74"
$ws.Range("E17").Value = "  +2.34%  "
$ws.Range("E18").Value = "  +0.47%  "
$ws.Range("D20").Value = "13.57"
$ws.Range("E20").Value = "  -0.54%  "
$ws.Range("D21").Value = "390.49"
$ws.Range("E21").Value = "  +1.59%  "
$ws.Range("D22").Value = "0.582"
$ws.Range("E22").Value = "  +3.03%  "
$ws.Range("D23").Value = "3.640.27"
$ws.Range("E23").Value = "  +0.35%  "
$ws.Range("D24").Value = "74.46"
$ws.Range("E24").Value = "  +0.43%  "
$ws.Range("E25").Value = "  -0.21%  "
$ws.Range("D26").Value = "5.66"
$ws.Range("E26").Value = "  -0.78%  "
$ws.Range("E27").Value = "  +2.72%  "
$ws.Range("E28").Value = "  +0.22%  "
$ws.Range("D29").Value = "7.42"
$ws.Range("E29").Value = "  +0.18%  "
$ws.Range("E30").Value = "  +1.92%  "
$ws.Range("E31").Value = "  -0.48%  "
$ws.Range("D32").Value = "1.48"
$ws.Range("E32").Value = "  -4.58%  "
$ws.Range("D33").Value = "0.157"
$ws.Range("E33").Value = "  +5.89%  "
$ws.Range("D34").Value = "3.528.46"
$ws.Range("E34").Value = "  +0.51%  "
$ws.Range("E36").Value = "  +0.20%  "
$ws.Range("D37").Value = "5.34"
$ws.Range("E37").Value = "  +1.71%  "
$ws.Range("D38").Value = "6.96"
$ws.Range("E38").Value = "  +1.77%  "
$ws.Range("E39").Value = "  +2.04%  "
$ws.Range("D40").Value = "165.82"
$ws.Range("E40").Value = "  +2.77%  "
$ws.Range("E41").Value = "  +1.82%  "
$ws.Range("D42").Value = "0.809"
$ws.Range("E42").Value = "  +0.66%  "
$ws.Range("E43").Value = "  +0.08%  "
$ws.Range("D44").Value = "4.45"
$ws.Range("E44").Value = "  +1.48%  "
$ws.Range("E45").Value = "  -2.34%  "
$ws.Range("D46").Value = "1.19"
$ws.Range("E46").Value = "  -0.39%  "
$ws.Range("E47").Value = "  +1.15%  "
$ws.Range("D48").Value = "0.927"
$ws.Range("E48").Value = "  +3.71%  "
$ws.Range("D49").Value = "2.417.35"
$ws.Range("E49").Value = "  -1.92%  "
$ws.Range("E50").Value = "  +1.34%  "
$ws.Range("E51").Value = "  +0.50%  "
